$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.178.17"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "3.640.87"
$ws.Range("E3").Value = "  +3.96%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'606.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "'202.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.220"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.19%  "
$ws.Range("D10").Value = "'0.650"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").Value = "'54.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "'0.0000307"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("D13").Value = "'9.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "4.225.89"
$ws.Range("E14").Value = "  +4.23%  "
$ws.Range("D15").Value = "'678.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +13.63%  "
$ws.Range("D16").Value = "71.282.19"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").Value = "'12.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").Value = "3.637.28"
$ws.Range("E18").Value = "  +4.53%  "
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").Value = "'18.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.60%  "
$ws.Range("D23").Value = "'5.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").Value = "'105.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("D25").Value = "'4.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("D27").Value = "'10.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("E28").Value = "  +5.52%  "
$ws.Range("D29").Value = "'34.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.27%  "
$ws.Range("D30").Value = "'4.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.27%  "
$ws.Range("D31").Value = "'7.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.95%  "
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").Value = "'63.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "0.0₃0878"
$ws.Range("E35").Value = "  +8.37%  "
$ws.Range("D36").Value = "3.920.87"
$ws.Range("E36").Value = "  +4.65%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'522.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.92%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  -5.37%  "
$ws.Range("D40").Value = "'3.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'0.392"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("E43").Value = "  +4.49%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0460"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.59%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'3.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.32%  "
$ws.Range("D46").Value = "'3.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.55%  "
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("D48").Value = "'8.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.75%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("E51").Value = "  +4.61%  "
